$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 999.4
$ws.Range("J19").Value = 999.5
$ws.Range("L19").Value = 999.5
$ws.Range("N19").Value = -1349.5
$ws.Range("H38").Value = 2074.3333
$ws.Range("I38").Value = 103.42857
$ws.Range("J38").Value = 8972.5
$ws.Range("K38").Value = 310.28571
$ws.Range("L38").Value = 26917.5
$ws.Range("M38").Value = 61.71429000000001
$ws.Range("N38").Value = -27661.5
$ws.Range("H40").Value = 15647761
$ws.Range("J40").Value = 33365372
$ws.Range("L40").Value = 33365372
$ws.Range("N40").Value = -33365722
$ws.Range("H42").Value = 125.71429
$ws.Range("I42").Value = 125.71429
$ws.Range("K42").Value = 377.14287
$ws.Range("M42").Value = -147.14287
$ws.Range("H64").Value = 20900044
$ws.Range("I64").Value = 7069349.5
$ws.Range("J64").Value = 35718644
$ws.Range("K64").Value = 7069349.5
$ws.Range("L64").Value = 35718644
$ws.Range("M64").Value = -7069101.5
$ws.Range("N64").Value = -35719140
$ws.Range("H67").Value = 20900044
$ws.Range("I67").Value = 7069349.5
$ws.Range("J67").Value = 35718644
$ws.Range("K67").Value = 7069349.5
$ws.Range("L67").Value = 35718644
$ws.Range("M67").Value = -7068491.5
$ws.Range("N67").Value = -35720360
$ws.Range("H98").Value = 2228.077
$ws.Range("I98").Value = 1543.1818
$ws.Range("J98").Value = 5995
$ws.Range("K98").Value = 1543.1818
$ws.Range("L98").Value = 5995
$ws.Range("M98").Value = -45.18180000000007
$ws.Range("N98").Value = -8991
$ws.Range("H121").Value = 4192.125
$ws.Range("J121").Value = 4192.125
$ws.Range("L121").Value = 12576.375
$ws.Range("N121").Value = -16070.375
$ws.Range("H122").Value = 2228.077
$ws.Range("I122").Value = 1543.1818
$ws.Range("J122").Value = 5995
$ws.Range("K122").Value = 4629.5454
$ws.Range("L122").Value = 17985
$ws.Range("M122").Value = -2179.5454
$ws.Range("N122").Value = -22885
$ws.Range("H132").Value = 20131.71
$ws.Range("I132").Value = 10740.9375
$ws.Range("J132").Value = 24305.389
$ws.Range("K132").Value = 32222.8125
$ws.Range("L132").Value = 72916.167
$ws.Range("M132").Value = -29692.8125
$ws.Range("N132").Value = -77976.167
$ws.Range("H137").Value = 10101454
$ws.Range("I137").Value = 527142.9399999999
$ws.Range("K137").Value = 1581428.82
$ws.Range("M137").Value = -1578878.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15365.371
$ws.Range("I32").Value = 15707.745
$ws.Range("K32").Value = 15707.745
$ws.Range("M32").Value = -15420.745
$ws.Range("H38").Value = 5000
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4533
$ws.Range("H107").Value = 67500
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 67500
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 67500
$ws.Range("M107").ClearContents() | Out-Null
$ws.Range("N107").Value = -75180
$ws.Range("H132").Value = 18386.584
$ws.Range("I132").Value = 23904.68
$ws.Range("K132").Value = 71714.04000000001
$ws.Range("M132").Value = -69184.04000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10248.5
$ws.Range("J7").Value = 8499
$ws.Range("L7").Value = 8499
$ws.Range("N7").Value = -8725
$ws.Range("H86").Value = 1716.1666
$ws.Range("I86").Value = 1942.1428
$ws.Range("J86").Value = 1399.8
$ws.Range("K86").Value = 1942.1428
$ws.Range("L86").Value = 1399.8
$ws.Range("M86").Value = -819.1428000000001
$ws.Range("N86").Value = -3645.8
$ws.Range("H89").Value = 1716.1666
$ws.Range("I89").Value = 1942.1428
$ws.Range("J89").Value = 1399.8
$ws.Range("K89").Value = 9710.714
$ws.Range("L89").Value = 6999
$ws.Range("M89").Value = -4094.714
$ws.Range("N89").Value = -18231
$ws.Range("H99").Value = 4167586.5
$ws.Range("I99").Value = 6945114
$ws.Range("J99").Value = 1295
$ws.Range("K99").Value = 6945114
$ws.Range("L99").Value = 1295
$ws.Range("M99").Value = -6943616
$ws.Range("N99").Value = -4291
$ws.Range("H134").Value = 3775.5293
$ws.Range("I134").Value = 1473.1818
$ws.Range("K134").Value = 4419.5454
$ws.Range("M134").Value = -1884.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3343332.2
$ws.Range("I6").Value = 10000000
$ws.Range("J6").Value = 14998.5
$ws.Range("K6").Value = 10000000
$ws.Range("L6").Value = 14998.5
$ws.Range("M6").Value = -9999887
$ws.Range("N6").Value = -15224.5
$ws.Range("H53").Value = 43791.5
$ws.Range("J53").Value = 43791.5
$ws.Range("L53").Value = 43791.5
$ws.Range("N53").Value = -45005.5
$ws.Range("H86").Value = 7483.3335
$ws.Range("I86").Value = 7400
$ws.Range("J86").Value = 7500
$ws.Range("K86").Value = 7400
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = -6277
$ws.Range("N86").Value = -9746
$ws.Range("H89").Value = 7483.3335
$ws.Range("I89").Value = 7400
$ws.Range("J89").Value = 7500
$ws.Range("K89").Value = 37000
$ws.Range("L89").Value = 37500
$ws.Range("M89").Value = -31384
$ws.Range("N89").Value = -48732
$ws.Range("H105").Value = 1624307.9
$ws.Range("I105").Value = 2066828.2
$ws.Range("K105").Value = 2066828.2
$ws.Range("M105").Value = -2065081.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 300
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents() | Out-Null
$ws.Range("H75").Value = 6264.909
$ws.Range("I75").Value = 2350
$ws.Range("J75").Value = 7134.8887
$ws.Range("K75").Value = 7050
$ws.Range("L75").Value = 21404.6661
$ws.Range("M75").Value = -6052
$ws.Range("N75").Value = -23400.6661
$ws.Range("H78").Value = 6264.909
$ws.Range("I78").Value = 2350
$ws.Range("J78").Value = 7134.8887
$ws.Range("K78").Value = 21150
$ws.Range("L78").Value = 64213.99830000001
$ws.Range("M78").Value = -16158
$ws.Range("N78").Value = -74197.99830000001
$ws.Range("H81").Value = 3815.6667
$ws.Range("I81").Value = 2614.3333
$ws.Range("J81").Value = 4416.3335
$ws.Range("K81").Value = 7842.999899999999
$ws.Range("L81").Value = 13249.0005
$ws.Range("M81").Value = -6719.999899999999
$ws.Range("N81").Value = -15495.0005
$ws.Range("H84").Value = 3815.6667
$ws.Range("I84").Value = 2614.3333
$ws.Range("J84").Value = 4416.3335
$ws.Range("K84").Value = 23528.9997
$ws.Range("L84").Value = 39747.0015
$ws.Range("M84").Value = -17912.9997
$ws.Range("N84").Value = -50979.0015
$ws.Range("H112").Value = 25000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 25000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 75000
$ws.Range("M112").ClearContents() | Out-Null
$ws.Range("N112").Value = -77216
$ws.Range("H121").Value = 1165.8
$ws.Range("J121").Value = 1157.75
$ws.Range("L121").Value = 3473.25
$ws.Range("N121").Value = -6093.25
$ws.Range("H134").Value = 12567.154
$ws.Range("I134").Value = 3714.2222
$ws.Range("K134").Value = 11142.6666
$ws.Range("M134").Value = -6072.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1327665.5
$ws.Range("I80").Value = 3127799.8
$ws.Range("K80").Value = 3127799.8
$ws.Range("M80").Value = -3126801.8
$ws.Range("H83").Value = 1327665.5
$ws.Range("I83").Value = 3127799.8
$ws.Range("K83").Value = 15638999
$ws.Range("M83").Value = -15634007
$ws.Range("H94").Value = 41862
$ws.Range("J94").Value = 34149.332
$ws.Range("L94").Value = 34149.332
$ws.Range("N94").Value = -35501.332
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents() | Out-Null
$ws.Range("H113").Value = 1654.5883
$ws.Range("I113").Value = 1339.4
$ws.Range("K113").Value = 1339.4
$ws.Range("M113").Value = 830.5999999999999
$ws.Range("H132").Value = 140237.27
$ws.Range("I132").Value = 254201.75
$ws.Range("K132").Value = 762605.25
$ws.Range("M132").Value = -760075.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 42858990
$ws.Range("J16").Value = 2832.889
$ws.Range("L16").Value = 2832.889
$ws.Range("N16").Value = -3172.889
$ws.Range("H61").Value = 1979.8334
$ws.Range("I61").Value = 1760.4736
$ws.Range("K61").Value = 1760.4736
$ws.Range("M61").Value = -1558.4736
$ws.Range("H113").Value = 1979.8334
$ws.Range("I113").Value = 1760.4736
$ws.Range("K113").Value = 1760.4736
$ws.Range("M113").Value = 409.5264
$ws.Range("H132").Value = 4657.5312
$ws.Range("I132").Value = 3574.8572
$ws.Range("J132").Value = 6724.4546
$ws.Range("K132").Value = 10724.5716
$ws.Range("L132").Value = 20173.3638
$ws.Range("M132").Value = -8194.571599999999
$ws.Range("N132").Value = -25233.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11880.952
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 11880.952
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 35642.856
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("N136").ClearContents() | Out-Null
